$wb = $excel.ActiveWorkbook

# Update the "Last Updated" timestamp on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 12:17 PM"

# Insert a new row for MIDWESTLTD at row 37 of the "Top Gainers" sheet,
# pushing the existing rows (HIRECT..GENESYS) down by one, and drop the
# row that falls off the bottom (old GENESYS row) to keep the same extent.
$gainers = $wb.Worksheets.Item("Top Gainers")
$gainers.Rows("37:37").Insert()
$gainers.Range("A37").Value = "🚀"
$gainers.Range("B37").Value = "MIDWESTLTD"
$gainers.Range("C37").Value = 3.8719
$gainers.Range("D37").Value = -0.7365
$gainers.Range("E37").Value = "N/A"
$gainers.Rows("77:77").Delete()
